$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")

# Deskcount value corrections
$ws.Range("C12").Value = 79
$ws.Range("C44").Value = 32
$ws.Range("C45").Value = 561

# "Include in Occupancy Calculation" flips from Yes to No for these locations
$ws.Range("F16").Value = "No"
$ws.Range("F22").Value = "No"
$ws.Range("F24").Value = "No"
$ws.Range("F37").Value = "No"
$ws.Range("F38").Value = "No"
$ws.Range("F47").Value = "No"
$ws.Range("F48").Value = "No"
$ws.Range("F49").Value = "No"

# Match the saved cursor/selection position from the authored workbook
$ws.Range("C42").Select()
